# Generate Report for Handback
# Update the handback status report: rows for
#   a501c29e-a94f-414d-b207-eec21fdbab64 (row 3)
#   cca6437e-08e6-4ef2-aecb-d9d34e532651 (row 5)
# got a fresh localization pass, so their timestamps move forward and the
# zh-cn priority flips from "ht" (human translation) to "mt" (machine
# translation).

$wb = $excel.ActiveWorkbook

# --- Overview sheet: "Latest HO Xliff Generate Date" column (G) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-08-23 14:17:16"
$wsOverview.Range("G5").Value = "2016-08-23 14:17:16"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H3").Value = "2016-08-23 14:16:59"
$wsZhCn.Range("H5").Value = "2016-08-23 14:16:59"
$wsZhCn.Range("K3").Value = "2016-08-23 14:17:33"
$wsZhCn.Range("K5").Value = "2016-08-23 14:17:33"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H3").Value = "2016-08-23 14:17:16"
$wsDeDe.Range("H5").Value = "2016-08-23 14:17:16"
$wsDeDe.Range("K3").Value = "2016-08-23 14:17:41"
$wsDeDe.Range("K5").Value = "2016-08-23 14:17:41"
